$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = "69.008.82"
$cell.Style = "Normal"
$ws.Range("E2").Value = "  +1.41%  "

$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = "3.752.83"
$cell.Style = "Normal"
$ws.Range("E3").Value = "  +2.15%  "

$ws.Range("E4").Value = "  +0.02%  "

$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "602.11"
$cell.Style = "Normal"
$ws.Range("E5").Value = "  +0.99%  "

$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "167.39"
$cell.Style = "Normal"
$ws.Range("E6").Value = "  +0.98%  "

$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = "3.750.68"
$cell.Style = "Normal"
$ws.Range("E7").Value = "  +2.17%  "

$ws.Range("E8").Value = "  -0.08%  "

$ws.Range("E9").Value = "  +1.52%  "

$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = "0.168"
$cell.Style = "Normal"
$ws.Range("E10").Value = "  +1.91%  "

$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = "6.46"
$cell.Style = "Normal"
$ws.Range("E11").Value = "  +3.31%  "

$ws.Range("E12").Value = "  +0.77%  "

$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = "38.02"
$cell.Style = "Normal"
$ws.Range("E13").Value = "  +0.32%  "

$ws.Range("E14").Value = "  +1.83%  "

$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = "4.381.43"
$cell.Style = "Normal"
$ws.Range("E15").Value = "  +2.17%  "

$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = "3.758.01"
$cell.Style = "Normal"
$ws.Range("E16").Value = "  +2.28%  "

$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = "69.052.82"
$cell.Style = "Normal"
$ws.Range("E17").Value = "  +1.48%  "

$ws.Range("E18").Value = "  +1.82%  "

$ws.Range("E19").Value = "  -0.95%  "

$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "17.24"
$cell.Style = "Normal"
$ws.Range("E20").Value = "  +0.84%  "

$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "10.90"
$cell.Style = "Normal"
$ws.Range("E21").Value = "  +20.45%  "

$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "492.48"
$cell.Style = "Normal"
$ws.Range("E22").Value = "  +0.79%  "

$ws.Range("E23").Value = "  +0.93%  "

$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "0.0000152"
$cell.Style = "Normal"
$ws.Range("E24").Value = "  +7.21%  "

$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "84.83"
$cell.Style = "Normal"
$ws.Range("E25").Value = "  +0.65%  "

$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = "2.30"
$cell.Style = "Normal"
$ws.Range("E26").Value = "  +0.93%  "

$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = "12.32"
$cell.Style = "Normal"
$ws.Range("E27").Value = "  +1.52%  "

$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = "10.12"
$cell.Style = "Normal"
$ws.Range("E28").Value = "  +1.02%  "

$ws.Range("E29").Value = "  +0.02%  "

$ws.Range("E30").Value = "  +3.74%  "

$ws.Range("E31").Value = "  +4.50%  "

$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = "8.01"
$cell.Style = "Normal"
$ws.Range("E32").Value = "  +2.54%  "

$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = "31.62"
$cell.Style = "Normal"
$ws.Range("E33").Value = "  +1.36%  "

$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = "3.897.22"
$cell.Style = "Normal"
$ws.Range("E34").Value = "  +2.19%  "

$ws.Range("B35").Value = "RenzoRestakedETH"
$ws.Range("C35").Value = "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = "3.687.47"
$cell.Style = "Normal"
$ws.Range("E35").Value = "  +1.98%  "

$ws.Range("B36").Value = "Hedera"
$ws.Range("C36").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = "0.109"
$cell.Style = "Normal"
$ws.Range("E36").Value = "  +0.62%  "

$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = "0.999"
$cell.Style = "Normal"
$ws.Range("E37").Value = "  -0.11%  "

$ws.Range("E38").Value = "  +2.07%  "

$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = "5.89"
$cell.Style = "Normal"
$ws.Range("E39").Value = "  +2.75%  "

$ws.Range("E40").Value = "  +2.15%  "

$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = "0.324"
$cell.Style = "Normal"
$ws.Range("E41").Value = "  +1.19%  "

$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = "2.95"
$cell.Style = "Normal"
$ws.Range("E42").Value = "  +5.03%  "

$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "429.83"
$cell.Style = "Normal"
$ws.Range("E43").Value = "  -0.59%  "

$ws.Range("B44").Value = "Stacks"
$ws.Range("C44").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "2.00"
$cell.Style = "Normal"
$ws.Range("E44").Value = "  +3.12%  "

$ws.Range("B45").Value = "OKB"
$ws.Range("C45").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = "48.65"
$cell.Style = "Normal"
$ws.Range("E45").Value = "  +0.37%  "

$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = "8.49"
$cell.Style = "Normal"
$ws.Range("E46").Value = "  +1.81%  "

$ws.Range("E47").Value = "  +0.00%  "

$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = "40.18"
$cell.Style = "Normal"
$ws.Range("E48").Value = "  -0.29%  "

$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = "141.36"
$cell.Style = "Normal"
$ws.Range("E49").Value = "  +0.00%  "

$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = "2.796.16"
$cell.Style = "Normal"
$ws.Range("E50").Value = "  +2.54%  "

$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "0.0352"
$cell.Style = "Normal"
$ws.Range("E51").Value = "  +1.62%  "
